$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '59.679.69'
Set-TextValue $ws.Range('E2') '  +6.27%  '

Set-TextValue $ws.Range('D3') '2.529.68'
Set-TextValue $ws.Range('E3') '  +5.56%  '

Set-TextValue $ws.Range('E4') '  +0.09%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue $ws.Range('D5') '506.13'
Set-TextValue $ws.Range('E5') '  +5.37%  '

$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws.Range('D6') '160.85'
Set-TextValue $ws.Range('E6') '  +7.89%  '

$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue $ws.Range('D7') '0.995'
Set-TextValue $ws.Range('E7') '  -0.30%  '

$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue $ws.Range('D8') '0.611'
Set-TextValue $ws.Range('E8') '  +22.29%  '

Set-TextValue $ws.Range('D9') '2.573.46'
Set-TextValue $ws.Range('E9') '  +7.15%  '

Set-TextValue $ws.Range('D10') '6.27'
Set-TextValue $ws.Range('E10') '  +14.85%  '

Set-TextValue $ws.Range('E11') '  +6.44%  '

Set-TextValue $ws.Range('D12') '0.342'
Set-TextValue $ws.Range('E12') '  +5.68%  '

Set-TextValue $ws.Range('D13') '0.126'
Set-TextValue $ws.Range('E13') '  +1.57%  '

Set-TextValue $ws.Range('D14') '2.964.97'
Set-TextValue $ws.Range('E14') '  +5.42%  '

Set-TextValue $ws.Range('D15') '59.422.25'
Set-TextValue $ws.Range('E15') '  +5.89%  '

Set-TextValue $ws.Range('D16') '22.02'
Set-TextValue $ws.Range('E16') '  +7.90%  '

Set-TextValue $ws.Range('D17') '0.0000139'
Set-TextValue $ws.Range('E17') '  +4.93%  '

Set-TextValue $ws.Range('D18') '2.552.08'
Set-TextValue $ws.Range('E18') '  +6.14%  '

Set-TextValue $ws.Range('D19') '4.76'
Set-TextValue $ws.Range('E19') '  +5.72%  '

Set-TextValue $ws.Range('D20') '333.65'
Set-TextValue $ws.Range('E20') '  +5.52%  '

Set-TextValue $ws.Range('D21') '10.30'
Set-TextValue $ws.Range('E21') '  +5.20%  '

Set-TextValue $ws.Range('D22') '6.06'
Set-TextValue $ws.Range('E22') '  +6.11%  '

Set-TextValue $ws.Range('E23') '  +0.40%  '

Set-TextValue $ws.Range('D24') '59.99'
Set-TextValue $ws.Range('E24') '  +5.57%  '

Set-TextValue $ws.Range('D25') '0.416'
Set-TextValue $ws.Range('E25') '  +5.16%  '

Set-TextValue $ws.Range('D26') '0.169'
Set-TextValue $ws.Range('E26') '  +6.68%  '

Set-TextValue $ws.Range('D27') '0.998'
Set-TextValue $ws.Range('E27') '  +0.27%  '

Set-TextValue $ws.Range('D28') '2.618.02'
Set-TextValue $ws.Range('E28') '  +4.60%  '

Set-TextValue $ws.Range('D29') '7.58'
Set-TextValue $ws.Range('E29') '  +3.99%  '

Set-TextValue $ws.Range('D30') '0.0₃0825'
Set-TextValue $ws.Range('E30') '  +6.54%  '

Set-TextValue $ws.Range('D31') '0.998'
Set-TextValue $ws.Range('E31') '  -0.16%  '

Set-TextValue $ws.Range('D32') '19.46'
Set-TextValue $ws.Range('E32') '  +8.10%  '

Set-TextValue $ws.Range('D33') '154.08'
Set-TextValue $ws.Range('E33') '  +3.78%  '

Set-TextValue $ws.Range('D34') '1.57'
Set-TextValue $ws.Range('E34') '  +5.51%  '

Set-TextValue $ws.Range('E35') '  +8.93%  '

Set-TextValue $ws.Range('D36') '3.94'
Set-TextValue $ws.Range('E36') '  +9.31%  '

Set-TextValue $ws.Range('E37') '  +7.61%  '

Set-TextValue $ws.Range('D38') '0.863'
Set-TextValue $ws.Range('E38') '  +2.22%  '

Set-TextValue $ws.Range('D39') '3.75'
Set-TextValue $ws.Range('E39') '  +10.59%  '

Set-TextValue $ws.Range('D40') '1.45'
Set-TextValue $ws.Range('E40') '  +6.86%  '

$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D41') '291.78'
Set-TextValue $ws.Range('E41') '  +13.60%  '

$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D42') '34.72'
Set-TextValue $ws.Range('E42') '  +3.73%  '

$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D43') '0.625'
Set-TextValue $ws.Range('E43') '  +6.83%  '

$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D44') '0.101'
Set-TextValue $ws.Range('E44') '  +7.04%  '

Set-TextValue $ws.Range('D45') '0.0558'
Set-TextValue $ws.Range('E45') '  +3.18%  '

Set-TextValue $ws.Range('D46') '0.995'
Set-TextValue $ws.Range('E46') '  -0.25%  '

Set-TextValue $ws.Range('E47') '  +6.70%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D48') '19.18'
Set-TextValue $ws.Range('E48') '  +12.52%  '

$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D49') '4.83'
Set-TextValue $ws.Range('E49') '  +3.47%  '

$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range('D50') '10.30'
Set-TextValue $ws.Range('E50') '  +0.82%  '

$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range('D51') '0.717'
Set-TextValue $ws.Range('E51') '  +12.53%  '
